$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "23.194.61"
$ws.Range("E2").Value = "  +0.47%  "

# Row 3 - Ethereum
Set-TextCell "D3" "1.602.47"
$ws.Range("E3").Value = "  -0.16%  "

# Row 4 - TetherUSD
Set-TextCell "D4" "1.000"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5 - USDC
$ws.Range("E5").Value = "  -0.14%  "

# Row 6 - BNB
Set-TextCell "D6" "303.38"
$ws.Range("E6").Value = "  +0.59%  "

# Row 7 - XRP
Set-TextCell "D7" "0.3784"
$ws.Range("E7").Value = "  +0.10%  "

# Row 8 - OKB
Set-TextCell "D8" "51.92"
$ws.Range("E8").Value = "  +3.10%  "

# Row 9 - Cardano
Set-TextCell "D9" "0.3626"
$ws.Range("E9").Value = "  -0.83%  "

# Row 10 - Polygon
$ws.Range("E10").Value = "  -0.29%  "

# Row 11 - BinanceUSD
Set-TextCell "D11" "1.000"
$ws.Range("E11").Value = "  -0.13%  "

# Row 12 - Dogecoin
Set-TextCell "D12" "0.08117"
$ws.Range("E12").Value = "  -0.34%  "

# Row 13 - Solana
Set-TextCell "D13" "22.83"
$ws.Range("E13").Value = "  -0.25%  "

# Row 14 - Polkadot
Set-TextCell "D14" "6.596"
$ws.Range("E14").Value = "  -0.22%  "

# Row 15 - Chainlink
$ws.Range("E15").Value = "  +0.26%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -1.15%  "

# Row 17 - WrappedEther
Set-TextCell "D17" "1.601.39"
$ws.Range("E17").Value = "  -0.03%  "

# Row 18 - Litecoin
Set-TextCell "D18" "93.87"
$ws.Range("E18").Value = "  +1.96%  "

# Row 19 - TRON
Set-TextCell "D19" "0.06876"
$ws.Range("E19").Value = "  +0.09%  "

# Row 20 - Avalanche
Set-TextCell "D20" "18.09"
$ws.Range("E20").Value = "  -1.13%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -0.75%  "

# Row 23 - Cosmos
$ws.Range("E23").Value = "  -1.12%  "

# Row 24 - WrappedBTC
Set-TextCell "D24" "23.187.89"
$ws.Range("E24").Value = "  +0.46%  "

# Row 25 - LidoDAOToken
Set-TextCell "D25" "3.019"
$ws.Range("E25").Value = "  +8.07%  "

# Row 26 - Toncoin
Set-TextCell "D26" "2.387"
$ws.Range("E26").Value = "  +1.38%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  +0.29%  "

# Row 28 - Monero
Set-TextCell "D28" "150.17"
$ws.Range("E28").Value = "  -0.27%  "

# Row 29 - HuobiToken
Set-TextCell "D29" "5.245"
$ws.Range("E29").Value = "  -0.64%  "

# Row 30 - BitcoinCash
Set-TextCell "D30" "133.99"
$ws.Range("E30").Value = "  +0.33%  "

# Row 31 - WEMIXTOKEN
Set-TextCell "D31" "2.371"
$ws.Range("E31").Value = "  -0.26%  "

# Row 32 - Filecoin
Set-TextCell "D32" "6.750"
$ws.Range("E32").Value = "  -1.25%  "

# Row 33 - WrappedliquidstakedEther2.0
Set-TextCell "D33" "1.779.09"
$ws.Range("E33").Value = "  +0.08%  "

# Row 34 - ImmutableX
Set-TextCell "D34" "0.9672"
$ws.Range("E34").Value = "  +1.11%  "

# Row 35 - Hedera
$ws.Range("E35").Value = "  -2.27%  "

# Row 36 - VeChain
Set-TextCell "D36" "0.02731"
$ws.Range("E36").Value = "  +0.07%  "

# Row 37 - FraxShare
Set-TextCell "D37" "10.22"
$ws.Range("E37").Value = "  -1.70%  "

# Row 38 - Algorand
Set-TextCell "D38" "0.2524"
$ws.Range("E38").Value = "  -1.09%  "

# Row 39 - Stellar
Set-TextCell "D39" "0.08799"
$ws.Range("E39").Value = "  -1.24%  "

# Row 40 - InternetComputer(DFINITY)
Set-TextCell "D40" "6.088"
$ws.Range("E40").Value = "  -2.95%  "

# Row 41 - now TrustWalletToken (was TheSandbox)
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell "D41" "1.365"
$ws.Range("E41").Value = "  -0.29%  "

# Row 42 - now TheSandbox (was TrustWalletToken)
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell "D42" "0.7106"
$ws.Range("E42").Value = "  +0.19%  "

# Row 43 - Aptos
Set-TextCell "D43" "12.57"
$ws.Range("E43").Value = "  -0.60%  "

# Row 44 - EnergySwap
Set-TextCell "D44" "15.65"
$ws.Range("E44").Value = "  +1.25%  "

# Row 45 - Decentraland
Set-TextCell "D45" "0.6554"
$ws.Range("E45").Value = "  -1.22%  "

# Row 46 - NEARProtocol
Set-TextCell "D46" "2.313"
$ws.Range("E46").Value = "  -0.55%  "

# Row 47 - PancakeSwap
Set-TextCell "D47" "4.019"
$ws.Range("E47").Value = "  +0.52%  "

# Row 48 - Quant
Set-TextCell "D48" "132.37"
$ws.Range("E48").Value = "  -0.22%  "

# Row 49 - Cronos
Set-TextCell "D49" "0.07949"
$ws.Range("E49").Value = "  +0.05%  "

# Row 50 - Flow
$ws.Range("E50").Value = "  -2.71%  "

# Row 51 - ThetaToken
Set-TextCell "D51" "1.211"
$ws.Range("E51").Value = "  +0.59%  "
